# Update DM integration test fixture
#
# 1. Bold the header row (row 1) on each of the three sheets.
# 2. Resize the columns on each sheet to their recalculated best-fit widths
#    (Excel widens the bold header text, which changes the auto-fit width).
# 3. Update the CodeSchemes ID (GUID) value in A2.

$wb = $excel.ActiveWorkbook

# --- CodeSchemes sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("CodeSchemes")
$ws1.Range("A1:N1").Font.Bold = $true

$ws1.Columns.Item(1).ColumnWidth  = 34.42857142857143
$ws1.Columns.Item(2).ColumnWidth  = 17.428571428571427
$ws1.Columns.Item(3).ColumnWidth  = 25.714285714285715
$ws1.Columns.Item(4).ColumnWidth  = 22.428571428571427
$ws1.Columns.Item(5).ColumnWidth  = 14.142857142857142
$ws1.Columns.Item(6).ColumnWidth  = 19.142857142857142
$ws1.Columns.Item(7).ColumnWidth  = 20.142857142857142
$ws1.Columns.Item(8).ColumnWidth  = 19.142857142857142
$ws1.Columns.Item(9).ColumnWidth  = 20.714285714285715
$ws1.Columns.Item(10).ColumnWidth = 24.0
$ws1.Columns.Item(11).ColumnWidth = 19.142857142857142
$ws1.Columns.Item(12).ColumnWidth = 15.857142857142858
$ws1.Columns.Item(13).ColumnWidth = 20.714285714285715
$ws1.Columns.Item(14).ColumnWidth = 27.285714285714285

$ws1.Range("A2").Value = "04867dfb-9eae-4665-8d46-f66b78260752"

# --- Codes sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Codes")
$ws2.Range("A1:J1").Font.Bold = $true

$ws2.Columns.Item(1).ColumnWidth  = 5.857142857142857
$ws2.Columns.Item(2).ColumnWidth  = 17.428571428571427
$ws2.Columns.Item(3).ColumnWidth  = 15.857142857142858
$ws2.Columns.Item(4).ColumnWidth  = 14.142857142857142
$ws2.Columns.Item(5).ColumnWidth  = 15.857142857142858
$ws2.Columns.Item(6).ColumnWidth  = 19.142857142857142
$ws2.Columns.Item(7).ColumnWidth  = 20.714285714285715
$ws2.Columns.Item(8).ColumnWidth  = 24.0
$ws2.Columns.Item(9).ColumnWidth  = 19.142857142857142
$ws2.Columns.Item(10).ColumnWidth = 15.857142857142858

# --- Extensions sheet ------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Extensions")
$ws3.Range("A1:I1").Font.Bold = $true

$ws3.Columns.Item(1).ColumnWidth = 5.857142857142857
$ws3.Columns.Item(2).ColumnWidth = 17.428571428571427
$ws3.Columns.Item(3).ColumnWidth = 14.142857142857142
$ws3.Columns.Item(4).ColumnWidth = 24.0
$ws3.Columns.Item(5).ColumnWidth = 15.857142857142858
$ws3.Columns.Item(6).ColumnWidth = 19.142857142857142
$ws3.Columns.Item(7).ColumnWidth = 19.142857142857142
$ws3.Columns.Item(8).ColumnWidth = 15.857142857142858
$ws3.Columns.Item(9).ColumnWidth = 24.0
